$wb = $excel.ActiveWorkbook

# --- Codebook sheet: update the "Favorite Typical Pizza Flavor" row ---
$codebook = $wb.Worksheets.Item("Codebook")
$codebook.Activate()

# Remove the leftover bold-looking (but actually non-bold) explicit style from
# the "Shoe Size" / "Favorite Typical Pizza Flavor" row labels so they fall
# back to the sheet's default formatting.
$codebook.Range("A5:A6").Style = "Normal"

# Update the variable definition text for the pizza-flavor row.
$codebook.Range("B6").Value = "favorite popular pizza flavors "

# Move the cursor/scroll position used before switching tabs.
$codebook.Range("A6").Select()

# --- Data sheet becomes the active tab again ---
$data = $wb.Worksheets.Item("Data")
$data.Activate()
$data.Range("E1").Select()
